$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 182 (everything from old row 182 downward shifts down by one)
$ws.Rows("182:182").Insert()

# Populate the newly inserted row 182 with the new "Merge Two Sorted Linked List" entry
$ws.Range("A182").Value() = "LinkedList"
$ws.Range("B182").Value() = "Merge Two Sorted Linked List"
$ws.Range("C182").Value() = "Yes"

# Match the row height used by the sibling "solved" row above it (row 181)
$ws.Rows("182:182").RowHeight = 19.5
